$wb = $excel.ActiveWorkbook

# Rename first sheet from "Sheet1" to "misc."
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "misc."

# Add header labels in K4 / L4, matching the style used by the other header cells (row 4)
$ws.Range("K4").Value = "other_indexes"
$ws.Range("L4").Value = "commodity"
$ws.Range("K4").Style = $ws.Range("J4").Style
$ws.Range("L4").Style = $ws.Range("J4").Style

# New data rows 11-12 (values entered in an order that reproduces the
# original shared-string table ordering: flo_emis, gas, co2captured, co2,
# "*ccs,*ccs-rf", "coal,oil")
$ws.Range("B11").Value = "flo_emis"
$ws.Range("D11").Value = "gas"
$ws.Range("L11").Value = "co2captured"
$ws.Range("K11").Value = "co2"
$ws.Range("E11").Value = "*ccs,*ccs-rf"
$ws.Range("H11").Value = 0.95

$ws.Range("B12").Value = "flo_emis"
$ws.Range("D12").Value = "coal,oil"
$ws.Range("L12").Value = "co2captured"
$ws.Range("K12").Value = "co2"
$ws.Range("E12").Value = "*ccs,*ccs-rf"
$ws.Range("H12").Value = 0.85

# Update selection to reflect new active cell
$ws.Range("D13").Select() | Out-Null

# Adjust column widths to match final layout (values chosen so the
# exported OOXML "width" lands on 10.33.. / 12 respectively)
$ws.Columns.Item(5).ColumnWidth = 9.5
$ws.Columns.Item(11).ColumnWidth = 11.166666666666666
